$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.946.44"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "1.787.31"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").Value = "226.14"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "30.24"
$ws.Range("E8").Value = "  -5.26%  "
$ws.Range("D9").Value = "46.70"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "0.0667"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "0.0924"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "2.045.31"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "1.788.10"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "0.625"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "10.39"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "34.001.52"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").Value = "4.19"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").Value = "69.09"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "252.27"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "0.0₃0740"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "10.33"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").Value = "4.22"
$ws.Range("E24").Value = "  -3.40%  "
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").Value = "158.35"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "16.51"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("D36").Value = "1.503.66"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "0.633"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "83.44"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "2.35"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "0.901"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").Value = "0.0520"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "2.04"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "1.941.58"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "5.73"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").Value = "  +9.07%  "
$ws.Range("D51").Value = "51.37"
$ws.Range("E51").Value = "  -5.91%  "
